$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.157.05"
$ws.Range("E2").Value = "  -1.58%  "

$ws.Range("D3").Value = "3.275.92"
$ws.Range("E3").Value = "  -1.71%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.38%  "

$ws.Range("E7").Value = "  +4.66%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -1.99%  "

$ws.Range("E10").Value = "  +1.29%  "

$ws.Range("E11").Value = "  -0.79%  "

$ws.Range("D12").Value = "3.847.49"
$ws.Range("E12").Value = "  -1.71%  "

$ws.Range("E13").Value = "  -3.56%  "

$ws.Range("D14").Value = "66.169.71"
$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.29%  "

$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "3.281.85"
$ws.Range("E17").Value = "  -1.54%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "434.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.21%  "

$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("E20").Value = "  -2.87%  "

$ws.Range("E21").Value = "  -3.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.08%  "

$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").Value = "3.420.40"
$ws.Range("E24").Value = "  -1.96%  "

$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("E27").Value = "  -5.45%  "

$ws.Range("E28").Value = "  -1.85%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("E30").Value = "  -1.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.60%  "

$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("E34").Value = "  -2.53%  "

$ws.Range("E35").Value = "  -3.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.70"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("E37").Value = "  -4.45%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "

$ws.Range("E39").Value = "  -2.84%  "

$ws.Range("D40").Value = "2.772.72"
$ws.Range("E40").Value = "  -2.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.773"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("E42").Value = "  -2.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.18%  "

$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "320.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("E47").Value = "  -2.75%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.90%  "

$ws.Range("E49").Value = "  -2.08%  "

$ws.Range("E50").Value = "  +3.19%  "

$ws.Range("E51").Value = "  +0.03%  "
